# separate pre mid and post mid timetables
# Applies the authored changes to sem1_DSAI_timetable_baskets.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "ELECTIVE_B1 [C304]"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "ELECTIVE_B1 [C304]"
$wsA.Range("E2").Value = "DS161 [C304]"
$wsA.Range("F2").Value = "EC161 [C201]"

$wsA.Range("B3").Value = "MA161 [C303]"
$wsA.Range("C3").Value = "MA162 [C002]"
$wsA.Range("D3").Value = "EC161 [C201]"
$wsA.Range("F3").Value = "DS161 [C304]"

$wsA.Range("B5").Value = "EC161 (Lab) [L306]"
$wsA.Range("C5").Value = "Free"
$wsA.Range("E5").Value = "MA161 [C303]"
$wsA.Range("F5").Value = "Free"

$wsA.Range("B6").Value = "EC161 (Lab) [L306]"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "Free"
$wsA.Range("D7").Value = "Free"
$wsA.Range("F7").Value = "MA162 [C002]"

# ---------------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "ELECTIVE_B1 [C201]"
$wsB.Range("D2").Value = "ELECTIVE_B1 [C201]"

$wsB.Range("D3").Value = "Free"

$wsB.Range("B5").Value = "MA161 [C205]"
$wsB.Range("C5").Value = "MA161 [C205]"
$wsB.Range("D5").Value = "EC161 (Lab) [L406]"
$wsB.Range("F5").Value = "Free"

$wsB.Range("D6").Value = "EC161 (Lab) [L406]"

$wsB.Range("B7").Value = "MA162 [C101]"
$wsB.Range("C7").Value = "DS161 [C403]"
$wsB.Range("D7").Value = "DS161 [C403]"
$wsB.Range("E7").Value = "Free"

$wsB.Range("C8").Value = "Free"

# ---------------------------------------------------------------------
# Verification_A (rows for MA162/EC161 swap places, rooms renamed)
# ---------------------------------------------------------------------
$wsVA = $wb.Worksheets.Item("Verification_A")

$wsVA.Range("I2").Value = "C304"
$wsVA.Range("I3").Value = "C303"

$wsVA.Range("A4").Value = "**EC161**"
$wsVA.Range("B4").Value = "Digital Design"
$wsVA.Range("C4").Value = "Prakash Pawar"
$wsVA.Range("D4").Value = "3-0-2-0-2"
$wsVA.Range("F4").Value = "2/1"
$wsVA.Range("H4").Value = "Partial"
$wsVA.Range("I4").Value = "L306, C201"

$wsVA.Range("A5").Value = "**MA162**"
$wsVA.Range("B5").Value = "Probability"
$wsVA.Range("C5").Value = "Chinmayananda"
$wsVA.Range("D5").Value = "2-0-0-0-2"
$wsVA.Range("F5").Value = "0/0"
$wsVA.Range("H5").Value = "Complete"
$wsVA.Range("I5").Value = "C002"

$wsVA.Range("I6").Value = "C304"

$wsVA.Range("H7").Value = "[WARN] 3 issues"

# ---------------------------------------------------------------------
# Verification_B (rows for DS161/MA161 swap places, rooms renamed)
# ---------------------------------------------------------------------
$wsVB = $wb.Worksheets.Item("Verification_B")

$wsVB.Range("I2").Value = "C201"
$wsVB.Range("I3").Value = "C402, L406"

$wsVB.Range("A4").Value = "**MA161**"
$wsVB.Range("B4").Value = "Statistics"
$wsVB.Range("C4").Value = "Ramesh Athe"
$wsVB.Range("D4").Value = "3-0-2-0-2"
$wsVB.Range("H4").Value = "Partial"
$wsVB.Range("I4").Value = "C205"

$wsVB.Range("A6").Value = "**DS161**"
$wsVB.Range("B6").Value = "Introduction to DATA science and artificial intelligence"
$wsVB.Range("C6").Value = "Girish Revadigar"
$wsVB.Range("D6").Value = "2-0-0-0-2"
$wsVB.Range("H6").Value = "Complete"
$wsVB.Range("I6").Value = "C403"

$wsVB.Range("H7").Value = "[WARN] 3 issues"

# ---------------------------------------------------------------------
# Room_Allocation (rooms re-sorted alphabetically, L407 row removed)
# ---------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Room_Allocation")

# Row 2: C002 (unchanged room id/capacity), reassigned to section A, course MA162
$wsR.Range("F2").Value = "A"
$wsR.Range("H2").Value = "MA162"

# Row 3: was C003 -> now C101
$wsR.Range("A3").Value = "C101"
$wsR.Range("B3").Value = "classroom"
$wsR.Range("C3").Value = "96"
$wsR.Range("F3").Value = "B"
$wsR.Range("H3").Value = "MA162"

# Row 4: was C101 -> now C201
$wsR.Range("A4").Value = "C201"
$wsR.Range("E4").Value = 4
$wsR.Range("F4").Value = "A, B"
$wsR.Range("G4").Value = 2
$wsR.Range("H4").Value = "ELECTIVE_B1, EC161"
$wsR.Range("I4").Value = "0.8"

# Row 5: was C202 -> now C205
$wsR.Range("A5").Value = "C205"
$wsR.Range("F5").Value = "B"
$wsR.Range("H5").Value = "MA161"

# Row 6: was C204 -> now C303
$wsR.Range("A6").Value = "C303"
$wsR.Range("F6").Value = "A"
$wsR.Range("H6").Value = "MA161"

# Row 7: was C205 -> now C304
$wsR.Range("A7").Value = "C304"
$wsR.Range("E7").Value = 4
$wsR.Range("F7").Value = "A"
$wsR.Range("G7").Value = 2
$wsR.Range("H7").Value = "DS161, ELECTIVE_B1"
$wsR.Range("I7").Value = "0.8"

# Row 8: was C303 -> now C402
$wsR.Range("A8").Value = "C402"
$wsR.Range("F8").Value = "B"
$wsR.Range("H8").Value = "EC161"

# Row 9: was C402 -> now C403
$wsR.Range("A9").Value = "C403"
$wsR.Range("C9").Value = "78"
$wsR.Range("E9").Value = 2
$wsR.Range("G9").Value = 1
$wsR.Range("H9").Value = "DS161"
$wsR.Range("I9").Value = "0.4"

# Row 10: was C403 -> now L306
$wsR.Range("A10").Value = "L306"
$wsR.Range("C10").Value = "96"
$wsR.Range("D10").Value = "Computers"
$wsR.Range("H10").Value = "EC161 (Lab)"

# Row 11: was L306 -> now L406
$wsR.Range("A11").Value = "L406"
$wsR.Range("C11").Value = "78"

# Row 12 (was L407) is no longer needed - remove it entirely
$wsR.Cells.Item(12, 1).EntireRow.Delete()

# ---------------------------------------------------------------------
# LTPSC_Compliance (status glyphs -> plain-text tags)
# ---------------------------------------------------------------------
$wsL = $wb.Worksheets.Item("LTPSC_Compliance")

$wsL.Range("G2").Value = "[OK]"
$wsL.Range("H2").Value = "[OK]"
$wsL.Range("I2").Value = "[OK]"
$wsL.Range("J2").Value = "[OK] FULLY COMPLIANT"

$wsL.Range("G3").Value = "[FAIL]"
$wsL.Range("H3").Value = "[OK]"
$wsL.Range("I3").Value = "[OK]"
$wsL.Range("J3").Value = "[WARN] PARTIAL"

$wsL.Range("G4").Value = "[FAIL]"
$wsL.Range("H4").Value = "[OK]"
$wsL.Range("I4").Value = "[OK]"
$wsL.Range("J4").Value = "[WARN] PARTIAL"

$wsL.Range("G5").Value = "[FAIL]"
$wsL.Range("H5").Value = "[OK]"
$wsL.Range("I5").Value = "[FAIL]"
$wsL.Range("J5").Value = "[WARN] PARTIAL"

$wsL.Range("G6").Value = "[OK]"
$wsL.Range("H6").Value = "[OK]"
$wsL.Range("I6").Value = "[OK]"
$wsL.Range("J6").Value = "[OK] FULLY COMPLIANT"

# ---------------------------------------------------------------------
# Executive_Summary
# ---------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Executive_Summary")

$wsE.Range("C3").Value = "2025-12-12 16:58"
$wsE.Range("C7").Value = "10/35"
$wsE.Range("D7").Value = "Utilization: 28.6%"
$wsE.Range("C9").Value = "[WARN] NEEDS REVIEW"
